# Updated cryptos list with refreshed Price (column D) and Volume(1h) (column E)
# values. Leading "'" forces Excel to store a numeric-looking Price as literal
# text (matching the workbook's existing text-based storage for these cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.383.99'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.848.09'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''240.23'
$ws.Range("D6").Value = '''0.6293'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.07630'
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value = '''0.2941'
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = '''24.46'
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("D11").Value = '''0.07745'
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.842.24'
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '''5.004'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = '''0.00001091'
$ws.Range("E14").Value = '  +9.10%  '
$ws.Range("D15").Value = '''0.6789'
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = '''83.42'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '2.091.01'
$ws.Range("E17").Value = '  -7.66%  '
$ws.Range("D18").Value = '''6.123'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '29.419.73'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''228.41'
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '''12.43'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D23").Value = '''7.444'
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '''157.23'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").Value = '''8.370'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = '''1.468'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '''1.296'
$ws.Range("E30").Value = '  +3.17%  '
$ws.Range("D31").Value = '''0.05621'
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = '''4.109'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '''4.041'
$ws.Range("D34").Value = '''1.849'
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").Value = '''1.155'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '''0.7089'
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").Value = '''2.588'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '''2.774'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '1.227.93'
$ws.Range("E39").Value = '  -2.04%  '
$ws.Range("D40").Value = '''0.01797'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("D41").Value = '''6.496'
$ws.Range("E41").Value = '  +4.59%  '
$ws.Range("D42").Value = '''0.9088'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '2.000.37'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = '''101.39'
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = '''66.00'
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("D48").Value = '''7.140'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").Value = '''0.4007'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("E51").Value = '  -0.20%  '
